$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume(1h)/Hora columns are stored as plain text (e.g. "303.60", "-0.14%", "17").
# Force a Text number format before writing so Excel does not silently turn the values
# into real numbers / percentages (which would also change "16" -> "17" style runs etc.).
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Column G ("Hora") moves from 16 to 17 for every data row.
$ws.Range("G2:G51").Value = "17"

# Updated Price (D) and Volume(1h) (E) figures per coin row.
$ws.Range("D2").Value = "303.60"
$ws.Range("E2").Value = "-0.14%"
$ws.Range("D3").Value = "37.15"
$ws.Range("E3").Value = "3.32%"
$ws.Range("D4").Value = "5.034"
$ws.Range("E4").Value = "-1.24%"
$ws.Range("D5").Value = "0.07844"
$ws.Range("E5").Value = "-0.20%"
$ws.Range("D6").Value = "2.214"
$ws.Range("E6").Value = "-3.18%"
$ws.Range("D7").Value = "8.001"
$ws.Range("E7").Value = "-1.00%"
$ws.Range("D8").Value = "4.029"
$ws.Range("E8").Value = "0.56%"
$ws.Range("D9").Value = "0.9284"
$ws.Range("E9").Value = "0.33%"
$ws.Range("D10").Value = "0.09821"
$ws.Range("E10").Value = "-2.33%"
$ws.Range("D11").Value = "0.1889"
$ws.Range("E11").Value = "3.17%"
$ws.Range("D12").Value = "0.08668"
$ws.Range("E12").Value = "0.15%"
$ws.Range("D13").Value = "0.03608"
$ws.Range("E13").Value = "6.04%"
$ws.Range("D14").Value = "0.09953"
$ws.Range("E14").Value = "0.50%"
$ws.Range("D15").Value = "0.001488"
$ws.Range("E15").Value = "1.36%"
$ws.Range("D16").Value = "0.005673"
$ws.Range("E16").Value = "1.43%"
$ws.Range("D17").Value = "3.462"
$ws.Range("E17").Value = "-0.86%"
$ws.Range("D18").Value = "2.297"
$ws.Range("E18").Value = "9.65%"
$ws.Range("D19").Value = "0.3433"
$ws.Range("E19").Value = "0.00%"
$ws.Range("D20").Value = "0.1326"
$ws.Range("E20").Value = "0.58%"
$ws.Range("D21").Value = "4.805"
$ws.Range("E21").Value = "5.22%"
$ws.Range("D22").Value = "0.2201"
$ws.Range("E22").Value = "-1.43%"
$ws.Range("D23").Value = "0.04594"
$ws.Range("E23").Value = "-1.38%"
$ws.Range("D24").Value = "0.005202"
$ws.Range("E24").Value = "15.83%"
$ws.Range("D26").Value = "0.0001402"
$ws.Range("E26").Value = "7.76%"
$ws.Range("D27").Value = "0.0002718"
$ws.Range("D39").Value = "0.01831"
$ws.Range("E39").Value = "4.43%"
$ws.Range("D40").Value = "0.04754"
$ws.Range("E40").Value = "1.31%"
$ws.Range("D41").Value = "0.007928"
$ws.Range("E41").Value = "0.71%"
$ws.Range("D42").Value = "0.1401"
$ws.Range("E42").Value = "-1.07%"
$ws.Range("D43").Value = "0.007566"
$ws.Range("E43").Value = "-13.90%"
$ws.Range("D44").Value = "0.002253"
$ws.Range("D45").Value = "0.01041"
$ws.Range("E45").Value = "13.13%"
$ws.Range("D46").Value = "0.00006335"
$ws.Range("E46").Value = "5.14%"
$ws.Range("E47").Value = "0.17%"
$ws.Range("D48").Value = "0.0005802"
$ws.Range("E48").Value = "0.02%"
$ws.Range("D49").Value = "35.64"
$ws.Range("E49").Value = "815.41%"
$ws.Range("D50").Value = "0.002690"
$ws.Range("D51").Value = "0.00002101"
$ws.Range("E51").Value = "0.17%"
